# New contact-form submission appended as the next row of the "Contacts"
# sheet (row 13), mirroring how the app's "generate Excel in memory" export
# keeps appending rows for each new submission.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

# Name
$ws.Cells.Item($row, 1).Value = "ahmed"

# Phone - force text storage (matches every other "Phone" cell in the
# sheet, which are stored as text even though they look numeric) by
# switching the cell to a text number format before typing the value,
# then reverting the format back to the sheet's normal style so no
# stray formatting is left behind.
$phoneCell = $ws.Cells.Item($row, 2)
$phoneCell.NumberFormat = "@"
$phoneCell.Value = "9854747474"
$phoneCell.NumberFormat = "General"
$phoneCell.Style = "Normal"

# Project
$ws.Cells.Item($row, 3).Value = "2BHK"

# Message - this submission left the message blank.
$ws.Cells.Item($row, 4).Value = ""

# Date
$ws.Cells.Item($row, 5).Value = "11/11/2025, 11:43:45 pm"
